$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 48.25514733333333
$ws.Range("H2").Value = 144.765442
$ws.Range("I2").Value = 0.9340796096783166
$ws.Range("J2").Value = 0.9340796096783167
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.45505566666667
$ws.Range("N2").Value = 31.365167
$ws.Range("O2").Value = 0.008202258778219224
$ws.Range("P2").Value = 0.008202258778219226
$ws.Range("Q2").Value = 504.5102515732016
$ws.Range("R2").Value = 4540.592264158814
$ws.Range("S2").Value = 0.007661562678039559
$ws.Range("T2").Value = 0.007661562678039561

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 48.25514733333333
$ws.Range("H3").Value = 144.765442
$ws.Range("I3").Value = 0.9340796096783166
$ws.Range("J3").Value = 0.9340796096783167
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 167.4277343333333
$ws.Range("N3").Value = 502.283203
$ws.Range("O3").Value = 0.1313513430666197
$ws.Range("P3").Value = 0.1313513430666197
$ws.Range("Q3").Value = 8079.249987941193
$ws.Range("R3").Value = 72713.24989147072
$ws.Range("S3").Value = 0.1226926112623908
$ws.Range("T3").Value = 0.1226926112623908

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 48.25514733333333
$ws.Range("H4").Value = 144.765442
$ws.Range("I4").Value = 0.9340796096783166
$ws.Range("J4").Value = 0.9340796096783167
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 227.495678
$ws.Range("N4").Value = 682.487034
$ws.Range("O4").Value = 0.1784761823728629
$ws.Range("P4").Value = 0.1784761823728629
$ws.Range("Q4").Value = 10977.83745958656
$ws.Range("R4").Value = 98800.53713627903
$ws.Range("S4").Value = 0.1667109627677199
$ws.Range("T4").Value = 0.1667109627677199

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 48.25514733333333
$ws.Range("H5").Value = 144.765442
$ws.Range("I5").Value = 0.9340796096783166
$ws.Range("J5").Value = 0.9340796096783167
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 147.9815903333333
$ws.Range("N5").Value = 443.9447709999999
$ws.Range("O5").Value = 0.1160953453549051
$ws.Range("P5").Value = 0.1160953453549051
$ws.Range("Q5").Value = 7140.873444155975
$ws.Range("R5").Value = 64267.86099740378
$ws.Range("S5").Value = 0.1084422948745791
$ws.Range("T5").Value = 0.1084422948745792

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 48.25514733333333
$ws.Range("H6").Value = 144.765442
$ws.Range("I6").Value = 0.9340796096783166
$ws.Range("J6").Value = 0.9340796096783167
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 665.7853396666667
$ws.Range("N6").Value = 1997.356019
$ws.Range("O6").Value = 0.5223256404173379
$ws.Range("P6").Value = 0.522325640417338
$ws.Range("Q6").Value = 32127.56965798838
$ws.Range("R6").Value = 289148.1269218954
$ws.Range("S6").Value = 0.4878937303260037
$ws.Range("T6").Value = 0.4878937303260039

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 48.25514733333333
$ws.Range("H7").Value = 144.765442
$ws.Range("I7").Value = 0.9340796096783166
$ws.Range("J7").Value = 0.9340796096783167
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 55.51027300000001
$ws.Range("N7").Value = 166.530819
$ws.Range("O7").Value = 0.043549230010055
$ws.Range("P7").Value = 0.04354923001005501
$ws.Range("Q7").Value = 2678.656402128556
$ws.Range("R7").Value = 24107.907619157
$ws.Range("S7").Value = 0.0406784477695834
$ws.Range("T7").Value = 0.04067844776958342

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.405489333333334
$ws.Range("H8").Value = 10.216468
$ws.Range("I8").Value = 0.06592039032168336
$ws.Range("J8").Value = 0.06592039032168336
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.45505566666667
$ws.Range("N8").Value = 31.365167
$ws.Range("O8").Value = 0.008202258778219224
$ws.Range("P8").Value = 0.008202258778219226
$ws.Range("Q8").Value = 35.60458055223956
$ws.Range("R8").Value = 320.441224970156
$ws.Range("S8").Value = 0.0005406961001796649
$ws.Range("T8").Value = 0.000540696100179665

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.405489333333334
$ws.Range("H9").Value = 10.216468
$ws.Range("I9").Value = 0.06592039032168336
$ws.Range("J9").Value = 0.06592039032168336
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 167.4277343333333
$ws.Range("N9").Value = 502.283203
$ws.Range("O9").Value = 0.1313513430666197
$ws.Range("P9").Value = 0.1313513430666197
$ws.Range("Q9").Value = 570.1733633763339
$ws.Range("R9").Value = 5131.560270387004
$ws.Range("S9").Value = 0.008658731804228907
$ws.Range("T9").Value = 0.008658731804228909

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.405489333333334
$ws.Range("H10").Value = 10.216468
$ws.Range("I10").Value = 0.06592039032168336
$ws.Range("J10").Value = 0.06592039032168336
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 227.495678
$ws.Range("N10").Value = 682.487034
$ws.Range("O10").Value = 0.1784761823728629
$ws.Range("P10").Value = 0.1784761823728629
$ws.Range("Q10").Value = 774.7341048084347
$ws.Range("R10").Value = 6972.606943275912
$ws.Range("S10").Value = 0.01176521960514307
$ws.Range("T10").Value = 0.01176521960514307

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.405489333333334
$ws.Range("H11").Value = 10.216468
$ws.Range("I11").Value = 0.06592039032168336
$ws.Range("J11").Value = 0.06592039032168336
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 147.9815903333333
$ws.Range("N11").Value = 443.9447709999999
$ws.Range("O11").Value = 0.1160953453549051
$ws.Range("P11").Value = 0.1160953453549051
$ws.Range("Q11").Value = 503.9497274098698
$ws.Range("R11").Value = 4535.547546688827
$ws.Range("S11").Value = 0.007653050480325972
$ws.Range("T11").Value = 0.007653050480325974

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.405489333333334
$ws.Range("H12").Value = 10.216468
$ws.Range("I12").Value = 0.06592039032168336
$ws.Range("J12").Value = 0.06592039032168336
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 665.7853396666667
$ws.Range("N12").Value = 1997.356019
$ws.Range("O12").Value = 0.5223256404173379
$ws.Range("P12").Value = 0.522325640417338
$ws.Range("Q12").Value = 2267.324872524544
$ws.Range("R12").Value = 20405.92385272089
$ws.Range("S12").Value = 0.03443191009133414
$ws.Range("T12").Value = 0.03443191009133414

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.405489333333334
$ws.Range("H13").Value = 10.216468
$ws.Range("I13").Value = 0.06592039032168336
$ws.Range("J13").Value = 0.06592039032168336
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 55.51027300000001
$ws.Range("N13").Value = 166.530819
$ws.Range("O13").Value = 0.043549230010055
$ws.Range("P13").Value = 0.04354923001005501
$ws.Range("Q13").Value = 189.0396425919214
$ws.Range("R13").Value = 1701.356783327292
$ws.Range("S13").Value = 0.002870782240471592
$ws.Range("T13").Value = 0.002870782240471592
